# The document has two logos that each appear twice (in the primary
# header/footer and in the first-page header/footer):
#   - the Pearson Edexcel logo (a .png) living in both footers, currently
#     named "image1.png", needs to become "image2.png"
#   - the BTEC logo (a .jpg) living in both headers, currently named
#     "image2.jpg", needs to become "image1.jpg"
#
# Walk every story range (main text + each header/footer story) and
# rename the inline picture found there based on which logo it is.
#
# Note: re-fetching the InlineShape through its own .Range.InlineShapes
# collection (rather than using the handle returned directly by the
# story's InlineShapes collection) is required for footer stories - Word
# occasionally reports the footer's first InlineShapes handle as stale,
# and reselecting it through its own Range freshens the reference before
# the rename is applied.

$d = $word.ActiveDocument

foreach ($story in $d.StoryRanges) {
    if ($story.InlineShapes.Count -eq 0) {
        continue
    }

    for ($i = 1; $i -le $story.InlineShapes.Count; $i++) {
        $pic = $story.InlineShapes.Item($i)
        $pic = $pic.Range.InlineShapes.Item(1)

        if ($pic.AlternativeText -like "*PearsonLogo.png") {
            $pic.Name = "image2.png"
        } elseif ($pic.AlternativeText -eq "BTec_Logo-Orange") {
            $pic.Name = "image1.jpg"
        }
    }
}
